# Applies the "Adds UX 1 and UX 2 results" commit to the UX 1.xlsx workbook.
#
# Summary of the real edit (from the xml diff):
#  - Sheet "RQ 1" (sheet2): A11 renamed from "Prototype 3 [ StackOverFlow UI ]"
#    to "Prototype 3 [ Tags ]"; row 17 gains a 6th answer (F17 = "P2"); a new
#    "Average" column (G) is added with an AVERAGE() formula for the three
#    rating rows (21-23), and a "Average" header label in G19.
#  - Sheet "RQ 2" (sheet3): same new "Average" column G with AVERAGE()
#    formulas for rows 7-9, and header label in G6.
#  - Sheet "RQ 3" (sheet4): same new "Average" column G with an AVERAGE()
#    formula for row 8, and header label in G7.
#  - The previously-removed shared string "Prototype 3 [ StackOverFlow UI ]"
#    disappears (superseded by "Prototype 3 [ Tags ]"), and "Average" is a
#    new shared string used 4 times.
#  - View state: "General" (sheet1) is no longer the active tab; "RQ 3"
#    (sheet4) becomes the active tab / selected sheet, with a couple of
#    incidental active-cell selection changes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "RQ 1"
# ---------------------------------------------------------------------
$rq1 = $wb.Worksheets.Item("RQ 1")

# Prototype 3 label was renamed.
$rq1.Cells.Item(11, 1).Value = "Prototype 3 [ Tags ]"

# Row 17 ("Which is convinient?") picks up a 6th ( U5 ) answer of "P2".
$rq1.Cells.Item(17, 6).Value = "P2"

# New "Average" column.
$rq1.Cells.Item(19, 7).Value = "Average"
$rq1.Range("G21").Formula = "=AVERAGE(B21:F21)"
$rq1.Range("G22").Formula = "=AVERAGE(B22:F22)"
$rq1.Range("G23").Formula = "=AVERAGE(B23:F23)"

# ---------------------------------------------------------------------
# Sheet "RQ 2"
# ---------------------------------------------------------------------
$rq2 = $wb.Worksheets.Item("RQ 2")

$rq2.Cells.Item(6, 7).Value = "Average"
$rq2.Range("G7").Formula = "=AVERAGE(B7:F7)"
$rq2.Range("G8").Formula = "=AVERAGE(B8:F8)"
$rq2.Range("G9").Formula = "=AVERAGE(B9:F9)"

# ---------------------------------------------------------------------
# Sheet "RQ 3"
# ---------------------------------------------------------------------
$rq3 = $wb.Worksheets.Item("RQ 3")

$rq3.Cells.Item(7, 7).Value = "Average"
$rq3.Range("G8").Formula = "=AVERAGE(B8:F8)"

# ---------------------------------------------------------------------
# View state: RQ 3 becomes the active sheet/tab, with updated selections
# on a few sheets (matches the final sheetView/selection state).
# ---------------------------------------------------------------------
$general = $wb.Worksheets.Item("General")
$general.Range("B15").Select()

$rq2.Range("H12").Select()

$rq3.Activate()
$rq3.Range("G12").Select()
